$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a cell to store a literal text value (no auto numeric conversion,
# no residual style index) matching the workbook's original inlineStr cells.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '26.517.04'
Set-TextValue $ws.Range("E2") '  -0.55%  '
Set-TextValue $ws.Range("D3") '1.627.54'
Set-TextValue $ws.Range("E3") '  -0.48%  '
Set-TextValue $ws.Range("E4") '  +0.19%  '
Set-TextValue $ws.Range("D5") '213.23'
Set-TextValue $ws.Range("E5") '  +0.12%  '
Set-TextValue $ws.Range("E6") '  +2.02%  '
Set-TextValue $ws.Range("E7") '  +0.18%  '
Set-TextValue $ws.Range("E8") '  -1.07%  '
Set-TextValue $ws.Range("D9") '0.0624'
Set-TextValue $ws.Range("E9") '  +0.01%  '
Set-TextValue $ws.Range("E10") '  -1.11%  '
Set-TextValue $ws.Range("E11") '  +0.53%  '
Set-TextValue $ws.Range("D12") '1.853.46'
Set-TextValue $ws.Range("E12") '  -0.54%  '
Set-TextValue $ws.Range("D13") '1.652.75'
Set-TextValue $ws.Range("E13") '  +0.75%  '
Set-TextValue $ws.Range("E14") '  +1.44%  '
Set-TextValue $ws.Range("D15") '0.524'
Set-TextValue $ws.Range("E15") '  -0.70%  '
Set-TextValue $ws.Range("D16") '65.10'
Set-TextValue $ws.Range("E16") '  +3.32%  '
Set-TextValue $ws.Range("D17") '26.536.60'
Set-TextValue $ws.Range("E17") '  -0.48%  '
Set-TextValue $ws.Range("D18") '0.0₃0741'
Set-TextValue $ws.Range("E18") '  +0.06%  '
Set-TextValue $ws.Range("D19") '214.74'
Set-TextValue $ws.Range("E19") '  +2.65%  '
Set-TextValue $ws.Range("E20") '  +0.19%  '
Set-TextValue $ws.Range("D21") '4.29'
Set-TextValue $ws.Range("E21") '  -0.59%  '
Set-TextValue $ws.Range("E22") '  +1.26%  '
Set-TextValue $ws.Range("D23") '9.31'
Set-TextValue $ws.Range("E23") '  -0.94%  '
Set-TextValue $ws.Range("D24") '2.12'
Set-TextValue $ws.Range("E24") '  +11.39%  '
Set-TextValue $ws.Range("D25") '147.73'
Set-TextValue $ws.Range("E25") '  +0.72%  '
Set-TextValue $ws.Range("E26") '  +0.21%  '
Set-TextValue $ws.Range("E27") '  -0.07%  '
Set-TextValue $ws.Range("D28") '6.90'
Set-TextValue $ws.Range("E28") '  +2.12%  '
Set-TextValue $ws.Range("D29") '15.55'
Set-TextValue $ws.Range("E29") '  +1.09%  '
Set-TextValue $ws.Range("D30") '0.0512'
Set-TextValue $ws.Range("E30") '  -1.53%  '
Set-TextValue $ws.Range("E31") '  -0.87%  '
Set-TextValue $ws.Range("E32") '  +3.41%  '
Set-TextValue $ws.Range("D34") '1.241.10'
Set-TextValue $ws.Range("E34") '  +6.13%  '
Set-TextValue $ws.Range("E35") '  -0.02%  '
Set-TextValue $ws.Range("E36") '  +0.49%  '
Set-TextValue $ws.Range("E37") '  +4.30%  '
Set-TextValue $ws.Range("E39") '  +0.90%  '
Set-TextValue $ws.Range("D40") '0.795'
Set-TextValue $ws.Range("E40") '  -1.61%  '
Set-TextValue $ws.Range("E41") '  -2.24%  '
Set-TextValue $ws.Range("E42") '  +0.48%  '
Set-TextValue $ws.Range("E43") '  -0.90%  '
Set-TextValue $ws.Range("D44") '1.763.31'
Set-TextValue $ws.Range("E44") '  -0.76%  '
Set-TextValue $ws.Range("D45") '93.19'
Set-TextValue $ws.Range("E45") '  +1.17%  '
Set-TextValue $ws.Range("D46") '1.59'
Set-TextValue $ws.Range("E46") '  +2.35%  '
Set-TextValue $ws.Range("B47") 'Aave'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D47") '54.89'
Set-TextValue $ws.Range("E47") '  +0.35%  '
Set-TextValue $ws.Range("B48") 'Cronos'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D48") '0.0510'
Set-TextValue $ws.Range("E48") '  -0.59%  '
Set-TextValue $ws.Range("B49") 'Algorand'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D49") '0.0958'
Set-TextValue $ws.Range("E49") '  +2.14%  '
Set-TextValue $ws.Range("B50") 'Mantle'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D50") '0.407'
Set-TextValue $ws.Range("E50") '  -0.63%  '
Set-TextValue $ws.Range("B51") 'EnergySwap'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D51") '7.49'
Set-TextValue $ws.Range("E51") '  -0.69%  '
